# Shift header values in row 1 (columns C:F) for every worksheet so that
# 'variable_trajectory_group' moves from F1 to C1, and the previous
# contents of C1, D1, E1 shift right by one column into D1, E1, F1
# respectively.
#
# Before: C1=normalize_group, D1=trajgroup_no_vary_q, E1=uniform_scaling_q, F1=variable_trajectory_group
# After:  C1=variable_trajectory_group, D1=normalize_group, E1=trajgroup_no_vary_q, F1=uniform_scaling_q

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "variable_trajectory_group"
    $ws.Range("D1").Value = "normalize_group"
    $ws.Range("E1").Value = "trajgroup_no_vary_q"
    $ws.Range("F1").Value = "uniform_scaling_q"
}
